# 2 best linear reg on data based on r value
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update column B with the new regression-selected data ---
# Row (A value) -> new B value, rows 2..30 (row 1 is unchanged)
$newB = @(
    1, 3, 2, 6, 5, 8, 7, 10, 13, 14, 15, 12, 13, 16, 17,
    14, 16, 20, 24, 32, 31, 29, 27, 34, 36, 47, 39, 46, 50
)
for ($i = 0; $i -lt $newB.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value2 = $newB[$i]
}

# --- Annotate the two best-fit rows (by r value) with a link back to their index ---
$ws.Hyperlinks.Add($ws.Range("C21"), "", "A20", "", "index 20")
$ws.Hyperlinks.Delete()
$ws.Range("C21").ClearFormats()
$ws.Range("C21").Value2 = "index 20"

$ws.Hyperlinks.Add($ws.Range("C30"), "", "A29", "", "index 29")
$ws.Hyperlinks.Delete()
$ws.Range("C30").ClearFormats()
$ws.Range("C30").Value2 = "index 29"

# --- Selection / view bookkeeping ---
[void]$ws.Range("B30").Select()

# --- Page setup to portrait orientation ---
$ws.PageSetup.Orientation = 1

# --- Window geometry (best effort; headless runtime may not persist this) ---
try {
    $win = $excel.ActiveWindow
    $win.Left = 5340
    $win.Top = 0
    $win.Width = 25140
    $win.Height = 15620
} catch {}
